$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("_requirements"); this shifts the
# previous columns B..J one position to the right (C..K).
$ws.Columns("B").Insert()

# Header + data for the new column B
$ws.Range("B1").Value = "_requirements"
$ws.Range("B5").Value = "l10n_it_reverse_charge"
$ws.Range("B6").Value = "l10n_it_split_payment"
$ws.Range("B7").Value = "l10n_it_dichiarazione_intento or l10n_it_lettera_intento"

# Widen the new column (target raw width 44.6 characters)
$ws.Columns("B").ColumnWidth = 43.83

# Match the final selection recorded in the workbook
$ws.Range("B8").Select()
